# Kilimanjaro_Weekly_Scoreboard.xlsx - append 4 new workout rows (302-305)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: copy the date formatting (style) from an existing date cell (B301)
# so the new date cells reuse the existing date number-format style instead
# of Excel creating a brand-new custom number format.
function Set-DateCell($addr, $serial) {
    $ws.Range("B301").Copy($ws.Range($addr))
    $ws.Range($addr).Value = $serial
}

# --- Row 302: Eric / Workout / Sauntering Hippo ---
$ws.Range("A302").Value = "Eric"
Set-DateCell "B302" 45501
$ws.Range("C302").Value = "Workout"
$ws.Range("D302").Value = 62
$ws.Range("E302").Value = 0
$ws.Range("F302").Value = 0
$ws.Range("G302").Value = 41
$ws.Range("H302").Value = 21
$ws.Range("I302").Value = 1
$ws.Range("J302").Value = 0
$ws.Range("K302").Value = 0
$ws.Range("L302").Value = "Sauntering Hippo"
$ws.Range("M302").Value = 7

# --- Row 303: Matt / Walk / Agile Antelope ---
$ws.Range("A303").Value = "Matt"
Set-DateCell "B303" 45501
$ws.Range("C303").Value = "Walk"
$ws.Range("D303").Value = 15
$ws.Range("E303").Value = 0.58
$ws.Range("F303").Value = 30
$ws.Range("G303").Value = 15
$ws.Range("H303").Value = 0
$ws.Range("I303").Value = 0
$ws.Range("J303").Value = 0
$ws.Range("K303").Value = 0
$ws.Range("L303").Value = "Agile Antelope"
$ws.Range("M303").Value = 7

# --- Row 304: Steven / Walk / Brave Leopard ---
$ws.Range("A304").Value = "Steven"
Set-DateCell "B304" 45502
$ws.Range("C304").Value = "Walk"
$ws.Range("D304").Value = 28
$ws.Range("E304").Value = 1.35
$ws.Range("F304").Value = 52
$ws.Range("G304").Value = 28
$ws.Range("H304").Value = 0
$ws.Range("I304").Value = 0
$ws.Range("J304").Value = 0
$ws.Range("K304").Value = 0
$ws.Range("L304").Value = "Brave Leopard"
$ws.Range("M304").Value = 8

# --- Row 305: Steven / Walk / Brave Leopard ---
$ws.Range("A305").Value = "Steven"
Set-DateCell "B305" 45503
$ws.Range("C305").Value = "Walk"
$ws.Range("D305").Value = 23
$ws.Range("E305").Value = 1.11
$ws.Range("F305").Value = 36
$ws.Range("G305").Value = 23
$ws.Range("H305").Value = 0
$ws.Range("I305").Value = 0
$ws.Range("J305").Value = 0
$ws.Range("K305").Value = 0
$ws.Range("L305").Value = "Brave Leopard"
$ws.Range("M305").Value = 8

# Update the frozen-pane scroll position and active selection to mirror
# the saved view state (scrolled down to show the newly added rows, with
# the final cell M305 selected).
$excel.Goto($ws.Range("A285"), $true) | Out-Null
$ws.Range("M305").Select() | Out-Null
